$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Selplg"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.156459
$ws.Range("H2").Value = 6.469377
$ws.Range("I2").Value = 0.01389412936885011
$ws.Range("J2").Value = 0.01392987523772938
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.112632333333332
$ws.Range("N2").Value = 27.337897
$ws.Range("O2").Value = 0.9981738658344552
$ws.Range("P2").Value = 0.9981738658344552
$ws.Range("Q2").Value = 19.65101800890766
$ws.Range("R2").Value = 176.859162080169
$ws.Range("S2").Value = 0.01386875682450915
$ws.Range("T2").Value = 0.01390443741663598

# Row 3: ECs -> M2
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Selplg"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.156459
$ws.Range("H3").Value = 6.469377
$ws.Range("I3").Value = 0.01389412936885011
$ws.Range("J3").Value = 0.01392987523772938
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01667133333333333
$ws.Range("N3").Value = 0.050014
$ws.Range("O3").Value = 0.001826134165544791
$ws.Range("P3").Value = 0.001826134165544791
$ws.Range("Q3").Value = 0.03595104680866666
$ws.Range("R3").Value = 0.323559421278
$ws.Range("S3").Value = 0.00002537254434095647
$ws.Range("T3").Value = 0.00002543782109339398

# Row 4: FAPs -> ECs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Selplg"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.009804666666666
$ws.Range("H4").Value = 9.029413999999999
$ws.Range("I4").Value = 0.01939226083762105
$ws.Range("J4").Value = 0.01944215192433629
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.112632333333332
$ws.Range("N4").Value = 27.337897
$ws.Range("O4").Value = 0.9981738658344552
$ws.Range("P4").Value = 0.9981738658344552
$ws.Range("Q4").Value = 27.42724332248422
$ws.Range("R4").Value = 246.845189902358
$ws.Range("S4").Value = 0.01935684796755831
$ws.Range("T4").Value = 0.01940664794645555

# Row 5: FAPs -> M2
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Selplg"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.009804666666666
$ws.Range("H5").Value = 9.029413999999999
$ws.Range("I5").Value = 0.01939226083762105
$ws.Range("J5").Value = 0.01944215192433629
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01667133333333333
$ws.Range("N5").Value = 0.050014
$ws.Range("O5").Value = 0.001826134165544791
$ws.Range("P5").Value = 0.001826134165544791
$ws.Range("Q5").Value = 0.05017745686622221
$ws.Range("R5").Value = 0.451597111796
$ws.Range("S5").Value = 0.00003541287006273604
$ws.Range("T5").Value = 0.00003550397788074291

# Row 6: M1 -> ECs
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Selplg"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 80.649016
$ws.Range("H6").Value = 241.947048
$ws.Range("I6").Value = 0.519624004803459
$ws.Range("J6").Value = 0.5209608580203196
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.112632333333332
$ws.Range("N6").Value = 27.337897
$ws.Range("O6").Value = 0.9981738658344552
$ws.Range("P6").Value = 0.9981738658344552
$ws.Range("Q6").Value = 734.9248308531172
$ws.Range("R6").Value = 6614.323477678055
$ws.Range("S6").Value = 0.5186751016550503
$ws.Range("T6").Value = 0.5200095135985772

# Row 7: M1 -> M2
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Selplg"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 80.649016
$ws.Range("H7").Value = 241.947048
$ws.Range("I7").Value = 0.519624004803459
$ws.Range("J7").Value = 0.5209608580203196
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01667133333333333
$ws.Range("N7").Value = 0.050014
$ws.Range("O7").Value = 0.001826134165544791
$ws.Range("P7").Value = 0.001826134165544791
$ws.Range("Q7").Value = 1.344526628741333
$ws.Range("R7").Value = 12.100739658672
$ws.Range("S7").Value = 0.0009489031484088073
$ws.Range("T7").Value = 0.0009513444217424346

# Row 8: M2 -> ECs
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Selplg"
$ws.Range("C8").Value = "Sele"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 68.19636666666668
$ws.Range("H8").Value = 204.5891
$ws.Range("I8").Value = 0.4393912153916231
$ws.Range("J8").Value = 0.4405216511573432
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.112632333333332
$ws.Range("N8").Value = 27.337897
$ws.Range("O8").Value = 0.9981738658344552
$ws.Range("P8").Value = 0.9981738658344552
$ws.Range("Q8").Value = 621.4484159025222
$ws.Range("R8").Value = 5593.0357431227
$ws.Range("S8").Value = 0.4385888280811562
$ws.Range("T8").Value = 0.4397171995195026

# Row 9: M2 -> M2
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Selplg"
$ws.Range("C9").Value = "Sele"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 68.19636666666668
$ws.Range("H9").Value = 204.5891
$ws.Range("I9").Value = 0.4393912153916231
$ws.Range("J9").Value = 0.4405216511573432
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01667133333333333
$ws.Range("N9").Value = 0.050014
$ws.Range("O9").Value = 0.001826134165544791
$ws.Range("P9").Value = 0.001826134165544791
$ws.Range("Q9").Value = 1.136924360822222
$ws.Range("R9").Value = 10.2323192474
$ws.Range("S9").Value = 0.0008023873104668933
$ws.Range("T9").Value = 0.0008044516378406285

# Row 10: sCs -> ECs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Selplg"
$ws.Range("C10").Value = "Sele"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.19484
$ws.Range("H10").Value = 2.38968
$ws.Range("I10").Value = 0.007698389598446743
$ws.Range("J10").Value = 0.005145463660271636
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.112632333333332
$ws.Range("N10").Value = 27.337897
$ws.Range("O10").Value = 0.9981738658344552
$ws.Range("P10").Value = 0.9981738658344552
$ws.Range("Q10").Value = 10.88813761716
$ws.Range("R10").Value = 65.32882570296
$ws.Range("S10").Value = 0.007684331306181345
$ws.Range("T10").Value = 0.005136067353284045

# Row 11: sCs -> M2
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Selplg"
$ws.Range("C11").Value = "Sele"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.19484
$ws.Range("H11").Value = 2.38968
$ws.Range("I11").Value = 0.007698389598446743
$ws.Range("J11").Value = 0.005145463660271636
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01667133333333333
$ws.Range("N11").Value = 0.050014
$ws.Range("O11").Value = 0.001826134165544791
$ws.Range("P11").Value = 0.001826134165544791
$ws.Range("Q11").Value = 0.01991957592
$ws.Range("R11").Value = 0.11951745552
$ws.Range("S11").Value = 0.00001405829226539825
$ws.Range("T11").Value = 0.00000939630698759119
